$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("instructions")
$wsData = $wb.Worksheets.Item("Data")
$wsCat = $wb.Worksheets.Item("Cat")

# --- "Data" sheet ---
# Fix the "Number of Observation(s)" label typo.
$wsData.Range("B2").Value = "Number of Observations"
$wsData.Range("I14").Select()

# --- "Cat" sheet ---
# Same label fix.
$wsCat.Range("B2").Value = "Number of Observations"
$wsCat.Range("F11").Select()

# --- "instructions" sheet ---
# Clear the stale threshold values in columns K (Improved Water Source),
# M (Improved Sanitation) and N (Access to handwashing facilities) for the
# data rows (3-7) -- these columns no longer carry threshold text.
$wsInstructions.Range("K3:K7").ClearContents()
$wsInstructions.Range("M3:M7").ClearContents()
$wsInstructions.Range("N3:N7").ClearContents()

# Update the active selection last so "instructions" remains the active
# (tabSelected) sheet, matching the saved workbook state.
$wsInstructions.Range("S2").Select()
